# Auto-generated script applying scheduled market-data refresh values
# to the Twintania_Profits leve-crafting profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 199.5
$ws.Range("I31").Value = 199.5
$ws.Range("K31").Value = 598.5
$ws.Range("M31").Value = -368.5

$ws.Range("H70").Value = 1262.75
$ws.Range("J70").Value = 1001.6667
$ws.Range("L70").Value = 3005.0001
$ws.Range("N70").Value = -3545.0001

$ws.Range("H73").Value = 1262.75
$ws.Range("J73").Value = 1001.6667
$ws.Range("L73").Value = 3005.0001
$ws.Range("N73").Value = -4877.0001

$ws.Range("H76").Value = 3594
$ws.Range("I76").Value = 3609.3333
$ws.Range("J76").Value = 3456
$ws.Range("K76").Value = 3609.3333
$ws.Range("L76").Value = 3456
$ws.Range("M76").Value = -3294.3333
$ws.Range("N76").Value = -4086

$ws.Range("H79").Value = 3594
$ws.Range("I79").Value = 3609.3333
$ws.Range("J79").Value = 3456
$ws.Range("K79").Value = 3609.3333
$ws.Range("L79").Value = 3456
$ws.Range("M79").Value = -2517.3333
$ws.Range("N79").Value = -5640

$ws.Range("H135").Value = 3501.3333
$ws.Range("I135").Value = 3783.4736
$ws.Range("K135").Value = 34051.2624
$ws.Range("M135").Value = -31516.2624

$ws.Range("H137").Value = 14788.833
$ws.Range("I137").Value = 6323.8184
$ws.Range("J137").Value = 21951.54
$ws.Range("K137").Value = 18971.4552
$ws.Range("L137").Value = 65854.62
$ws.Range("M137").Value = -16421.4552
$ws.Range("N137").Value = -70954.62

$ws.Range("H141").Value = 2849.4707
$ws.Range("I141").Value = 2770.6365
$ws.Range("J141").Value = 2994
$ws.Range("K141").Value = 8311.9095
$ws.Range("L141").Value = 8982
$ws.Range("M141").Value = -3131.9095
$ws.Range("N141").Value = -19342

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3485.5193
$ws.Range("I2").Value = 3782.3438
$ws.Range("K2").Value = 3782.3438
$ws.Range("M2").Value = -3669.3438

$ws.Range("H45").Value = 7399.75
$ws.Range("I45").Value = 8920.066
$ws.Range("J45").Value = 2838.8
$ws.Range("K45").Value = 8920.066
$ws.Range("L45").Value = 2838.8
$ws.Range("M45").Value = -8543.066
$ws.Range("N45").Value = -3592.8

$ws.Range("H116").Value = 3485.5193
$ws.Range("I116").Value = 3782.3438
$ws.Range("K116").Value = 3782.3438
$ws.Range("M116").Value = -1488.3438

$ws.Range("H132").Value = 8937.167
$ws.Range("I132").Value = 5124.6
$ws.Range("J132").Value = 28000
$ws.Range("K132").Value = 15373.8
$ws.Range("L132").Value = 84000
$ws.Range("M132").Value = -12843.8
$ws.Range("N132").Value = -89060

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3485.5193
$ws.Range("I3").Value = 3782.3438
$ws.Range("K3").Value = 3782.3438
$ws.Range("M3").Value = -3668.3438

$ws.Range("H105").Value = 3639.3794
$ws.Range("I105").Value = 3223.2173
$ws.Range("K105").Value = 3223.2173
$ws.Range("M105").Value = -1476.2173

$ws.Range("H107").Value = 4266.3335
$ws.Range("I107").Value = 3899.5
$ws.Range("J107").Value = 5000
$ws.Range("K107").Value = 3899.5
$ws.Range("L107").Value = 5000
$ws.Range("M107").Value = -1979.5
$ws.Range("N107").Value = -8840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 592.6667
$ws.Range("I22").Value = 188.23077
$ws.Range("K22").Value = 188.23077
$ws.Range("M22").Value = 161.76923

$ws.Range("H23").Value = 30873.75
$ws.Range("I23").Value = 3999.5
$ws.Range("J23").Value = 39831.832
$ws.Range("K23").Value = 3999.5
$ws.Range("L23").Value = 39831.832
$ws.Range("M23").Value = -3759.5
$ws.Range("N23").Value = -40311.832

$ws.Range("H27").Value = 30873.75
$ws.Range("I27").Value = 3999.5
$ws.Range("J27").Value = 39831.832
$ws.Range("K27").Value = 3999.5
$ws.Range("L27").Value = 39831.832
$ws.Range("M27").Value = -3807.5
$ws.Range("N27").Value = -40215.832

$ws.Range("H52").Value = 67499.75
$ws.Range("I52").Value = 39999
$ws.Range("J52").Value = 76666.664
$ws.Range("K52").Value = 39999
$ws.Range("L52").Value = 76666.664
$ws.Range("M52").Value = -39705
$ws.Range("N52").Value = -77254.664

$ws.Range("H58").Value = 2868.0469
$ws.Range("I58").Value = 1985.8718
$ws.Range("J58").Value = 4244.24
$ws.Range("K58").Value = 1985.8718
$ws.Range("L58").Value = 4244.24
$ws.Range("M58").Value = -1782.8718
$ws.Range("N58").Value = -4650.24

$ws.Range("H86").Value = 3129.8333
$ws.Range("I86").Value = 2442.1538
$ws.Range("K86").Value = 2442.1538
$ws.Range("M86").Value = -1319.1538

$ws.Range("H89").Value = 3129.8333
$ws.Range("I89").Value = 2442.1538
$ws.Range("K89").Value = 12210.769
$ws.Range("M89").Value = -6594.769

$ws.Range("H94").Value = 3595.3333
$ws.Range("I94").Value = 3595.3333
$ws.Range("K94").Value = 3595.3333
$ws.Range("M94").Value = -3144.3333

$ws.Range("H105").Value = 1824.0834
$ws.Range("I105").Value = 1893.9
$ws.Range("K105").Value = 1893.9
$ws.Range("M105").Value = -146.9000000000001

$ws.Range("H107").Value = 1016.2857
$ws.Range("I107").Value = 1009.8461
$ws.Range("J107").Value = 1100
$ws.Range("K107").Value = 1009.8461
$ws.Range("L107").Value = 1100
$ws.Range("M107").Value = 910.1539
$ws.Range("N107").Value = -4940

$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("M135").ClearContents()

$ws.Range("H136").Value = 2868.0469
$ws.Range("I136").Value = 1985.8718
$ws.Range("J136").Value = 4244.24
$ws.Range("K136").Value = 5957.6154
$ws.Range("L136").Value = 12732.72
$ws.Range("M136").Value = -3407.6154
$ws.Range("N136").Value = -17832.72

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1430.5333
$ws.Range("J5").Value = 2225.5
$ws.Range("L5").Value = 6676.5
$ws.Range("N5").Value = -6900.5

$ws.Range("H88").Value = 8559.4
$ws.Range("I88").Value = 6965.6665
$ws.Range("K88").Value = 20896.9995
$ws.Range("M88").Value = -20468.9995

$ws.Range("H91").Value = 8559.4
$ws.Range("I91").Value = 6965.6665
$ws.Range("K91").Value = 20896.9995
$ws.Range("M91").Value = -19414.9995

$ws.Range("H98").Value = 1366.1
$ws.Range("I98").Value = 1439.7142
$ws.Range("J98").Value = 1326.4615
$ws.Range("K98").Value = 4319.142599999999
$ws.Range("L98").Value = 3979.3845
$ws.Range("M98").Value = -2821.142599999999
$ws.Range("N98").Value = -6975.3845

$ws.Range("H117").Value = 2703.7144
$ws.Range("I117").Value = 1385.2
$ws.Range("J117").Value = 6000
$ws.Range("K117").Value = 4155.6
$ws.Range("L117").Value = 18000
$ws.Range("M117").Value = -713.6000000000004
$ws.Range("N117").Value = -24884

$ws.Range("H121").Value = 3151.6785
$ws.Range("J121").Value = 3565.5833
$ws.Range("L121").Value = 10696.7499
$ws.Range("N121").Value = -13316.7499

$ws.Range("H135").Value = 1430.5333
$ws.Range("J135").Value = 2225.5
$ws.Range("L135").Value = 20029.5
$ws.Range("N135").Value = -25099.5

$ws.Range("H140").Value = 989.5714
$ws.Range("I140").Value = 1017.8333
$ws.Range("J140").Value = 820
$ws.Range("K140").Value = 3053.4999
$ws.Range("L140").Value = 2460
$ws.Range("M140").Value = 2126.5001
$ws.Range("N140").Value = -12820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 49994
$ws.Range("J52").Value = 49994
$ws.Range("L52").Value = 49994
$ws.Range("N52").Value = -50512

$ws.Range("H93").Value = 42997.668
$ws.Range("J93").Value = 42997.668
$ws.Range("L93").Value = 42997.668
$ws.Range("N93").Value = -46741.668

$ws.Range("H113").Value = 504048.25
$ws.Range("I113").Value = 504048.25
$ws.Range("K113").Value = 504048.25
$ws.Range("M113").Value = -501878.25

$ws.Range("H132").Value = 39982
$ws.Range("I132").Value = 45178.4
$ws.Range("K132").Value = 135535.2
$ws.Range("M132").Value = -133005.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9599.8
$ws.Range("I7").Value = 7999.6665
$ws.Range("J7").Value = 12000
$ws.Range("K7").Value = 7999.6665
$ws.Range("L7").Value = 12000
$ws.Range("M7").Value = -7887.6665
$ws.Range("N7").Value = -12224

$ws.Range("H126").Value = 9599.8
$ws.Range("I126").Value = 7999.6665
$ws.Range("J126").Value = 12000
$ws.Range("K126").Value = 23998.9995
$ws.Range("L126").Value = 36000
$ws.Range("M126").Value = -21528.9995
$ws.Range("N126").Value = -40940

$ws.Range("H132").Value = 5664.2856
$ws.Range("I132").Value = 5496.25
$ws.Range("J132").Value = 6672.5
$ws.Range("K132").Value = 16488.75
$ws.Range("L132").Value = 20017.5
$ws.Range("M132").Value = -13958.75
$ws.Range("N132").Value = -25077.5

$ws.Range("H136").Value = 5188.067
$ws.Range("I136").Value = 4896.1787
$ws.Range("K136").Value = 14688.5361
$ws.Range("M136").Value = -12138.5361

$ws.Range("H140").Value = 78880.164
$ws.Range("J140").Value = 78880.164
$ws.Range("L140").Value = 78880.164
$ws.Range("N140").Value = -89240.164

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10004
$ws.Range("J3").Value = 10004
$ws.Range("L3").Value = 10004
$ws.Range("N3").Value = -10232

$ws.Range("H42").Value = 74998
$ws.Range("I42").Value = 49997
$ws.Range("K42").Value = 49997
$ws.Range("M42").Value = -49619

$ws.Range("H43").Value = 55124.75
$ws.Range("I43").Value = 40166.668
$ws.Range("K43").Value = 40166.668
$ws.Range("M43").Value = -40017.668

$ws.Range("H132").Value = 158635.36
$ws.Range("I132").Value = 227867.78
$ws.Range("J132").Value = 31247.68
$ws.Range("K132").Value = 683603.34
$ws.Range("L132").Value = 93743.04000000001
$ws.Range("M132").Value = -681073.34
$ws.Range("N132").Value = -98803.04000000001

$ws.Range("H136").Value = 7144622.5
$ws.Range("I136").Value = 8697111
$ws.Range("J136").Value = 3175.2
$ws.Range("K136").Value = 26091333
$ws.Range("L136").Value = 9525.599999999999
$ws.Range("M136").Value = -26088783
$ws.Range("N136").Value = -14625.6
